$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new "Formula Text" column (reuse B1's bold/centered formatting)
$ws.Range("D1").Value = "Formula Text"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108

# FORMULATEXT formulas showing the text of the B-column formulas.
# D2 is entered on its own; D3:D7 are entered together as a single range so
# Excel stores them as one shared formula group (mirrors B4:B6 / B5 groups).
$ws.Range("D2").Formula = "=FORMULATEXT(B2)"
$ws.Range("D3:D7").Formula = "=FORMULATEXT(B3)"

# Move the active selection to B7
$ws.Range("B7").Select()
